# Generate Report for Handoff
# The file "b76551b3-5742-46c8-b8a3-f640ef6031e6.md" has moved from "In Translation"
# to "Ready for handoff" for both the zh-cn and de-de locales, and its handoff
# xliff files + datetimes were refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the b76551b3-... file. Update its per-locale status
# (columns E = zh-cn, F = de-de) and the "Latest HO Xliff Generate Date" (G).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 04:05:46"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is the b76551b3-... file.
#   C = Status, E = Priority, H = Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2017-01-03 04:05:37"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is the b76551b3-... file.
#   C = Status, E = Priority, H = Latest Handoff Datetime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2017-01-03 04:05:46"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
